$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 173.5
$ws.Range("I33").Value = 148.15384
$ws.Range("K33").Value = 148.15384
$ws.Range("M33").Value = 80.84616
$ws.Range("H55").Value = 325.30768
$ws.Range("I55").Value = 333.2
$ws.Range("K55").Value = 333.2
$ws.Range("M55").Value = -119.2
$ws.Range("H70").Value = 6659.091
$ws.Range("I70").Value = 6653.846
$ws.Range("K70").Value = 19961.538
$ws.Range("M70").Value = -19691.538
$ws.Range("H73").Value = 6659.091
$ws.Range("I73").Value = 6653.846
$ws.Range("K73").Value = 19961.538
$ws.Range("M73").Value = -19025.538
$ws.Range("H98").Value = 1243.7222
$ws.Range("I98").Value = 1212.8
$ws.Range("K98").Value = 1212.8
$ws.Range("M98").Value = 285.2
$ws.Range("H106").Value = 29415.666
$ws.Range("I106").Value = 26159.643
$ws.Range("K106").Value = 26159.643
$ws.Range("M106").Value = -25528.643
$ws.Range("H122").Value = 1243.7222
$ws.Range("I122").Value = 1212.8
$ws.Range("K122").Value = 3638.4
$ws.Range("M122").Value = -1188.4
$ws.Range("H138").Value = 1999.2363
$ws.Range("I138").Value = 1701.0714
$ws.Range("K138").Value = 5103.2142
$ws.Range("M138").Value = 36.78579999999965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 939.8
$ws.Range("I74").Value = 950
$ws.Range("K74").Value = 950
$ws.Range("M74").Value = -76
$ws.Range("H77").Value = 939.8
$ws.Range("I77").Value = 950
$ws.Range("K77").Value = 4750
$ws.Range("M77").Value = -382
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("M122").Value = -3548.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 808
$ws.Range("I29").Value = 808
$ws.Range("K29").Value = 808
$ws.Range("M29").Value = -519
$ws.Range("H105").Value = 2852.6333
$ws.Range("I105").Value = 2267.4092
$ws.Range("K105").Value = 2267.4092
$ws.Range("M105").Value = -520.4092000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 5000
$ws.Range("J57").Value = 5000
$ws.Range("L57").Value = 5000
$ws.Range("N57").Value = -6120
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H132").Value = 8000
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 114866.336
$ws.Range("I23").Value = 2999.5
$ws.Range("J23").Value = 146828.28
$ws.Range("K23").Value = 8998.5
$ws.Range("L23").Value = 440484.84
$ws.Range("M23").Value = -8763.5
$ws.Range("N23").Value = -440954.84
$ws.Range("H75").Value = 1211.3334
$ws.Range("I75").Value = 1200
$ws.Range("J75").Value = 1217
$ws.Range("K75").Value = 3600
$ws.Range("L75").Value = 3651
$ws.Range("M75").Value = -2602
$ws.Range("N75").Value = -5647
$ws.Range("H78").Value = 1211.3334
$ws.Range("I78").Value = 1200
$ws.Range("J78").Value = 1217
$ws.Range("K78").Value = 10800
$ws.Range("L78").Value = 10953
$ws.Range("M78").Value = -5808
$ws.Range("N78").Value = -20937
$ws.Range("H114").Value = 1000
$ws.Range("I114").Value = 1000
$ws.Range("K114").Value = 3000
$ws.Range("M114").Value = 254
$ws.Range("H117").Value = 1287.4546
$ws.Range("J117").Value = 1526.375
$ws.Range("L117").Value = 4579.125
$ws.Range("N117").Value = -11463.125
$ws.Range("H137").Value = 3977.4443
$ws.Range("J137").Value = 3449.75
$ws.Range("L137").Value = 10349.25
$ws.Range("N137").Value = -20549.25
$ws.Range("H139").Value = 1631.5
$ws.Range("I139").Value = 30
$ws.Range("K139").Value = 90
$ws.Range("M139").Value = 5050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 15411.8
$ws.Range("I26").Value = 8000
$ws.Range("K26").Value = 8000
$ws.Range("M26").Value = -7720
$ws.Range("H50").Value = 15411.8
$ws.Range("I50").Value = 8000
$ws.Range("K50").Value = 8000
$ws.Range("M50").Value = -7502
$ws.Range("H132").Value = 2745.6453
$ws.Range("I132").Value = 3155.625
$ws.Range("K132").Value = 9466.875
$ws.Range("M132").Value = -6936.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3752
$ws.Range("I7").Value = 5499.5
$ws.Range("J7").Value = 2004.5
$ws.Range("K7").Value = 5499.5
$ws.Range("L7").Value = 2004.5
$ws.Range("M7").Value = -5387.5
$ws.Range("N7").Value = -2228.5
$ws.Range("H9").Value = 345
$ws.Range("I9").Value = 126.666664
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 126.666664
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 97.333336
$ws.Range("N9").Value = -1448
$ws.Range("H22").Value = 2602.318
$ws.Range("I22").Value = 1575.0834
$ws.Range("J22").Value = 3835
$ws.Range("K22").Value = 1575.0834
$ws.Range("L22").Value = 3835
$ws.Range("M22").Value = -1280.0834
$ws.Range("N22").Value = -4425
$ws.Range("H27").Value = 2602.318
$ws.Range("I27").Value = 1575.0834
$ws.Range("J27").Value = 3835
$ws.Range("K27").Value = 1575.0834
$ws.Range("L27").Value = 3835
$ws.Range("M27").Value = -1468.0834
$ws.Range("N27").Value = -4049
$ws.Range("H35").Value = 5333.3335
$ws.Range("I35").Value = 500
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 500
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = -164
$ws.Range("N35").Value = -15672
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29593
$ws.Range("H82").Value = 336330.66
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 336330.66
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 336330.66
$ws.Range("N82").Value = -337052.66
$ws.Range("H85").Value = 336330.66
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 336330.66
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 336330.66
$ws.Range("N85").Value = -338826.66
$ws.Range("H93").Value = 567.3333
$ws.Range("I93").Value = 351
$ws.Range("K93").Value = 351
$ws.Range("M93").Value = 897
$ws.Range("H122").Value = 5469.591
$ws.Range("I122").Value = 2866.6667
$ws.Range("K122").Value = 8600.000100000001
$ws.Range("M122").Value = -6150.000100000001
$ws.Range("H126").Value = 3752
$ws.Range("I126").Value = 5499.5
$ws.Range("J126").Value = 2004.5
$ws.Range("K126").Value = 16498.5
$ws.Range("L126").Value = 6013.5
$ws.Range("M126").Value = -14028.5
$ws.Range("N126").Value = -10953.5
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 30785.25
$ws.Range("J33").Value = 31380.334
$ws.Range("L33").Value = 31380.334
$ws.Range("N33").Value = -31880.334
$ws.Range("H36").Value = 30785.25
$ws.Range("J36").Value = 31380.334
$ws.Range("L36").Value = 31380.334
$ws.Range("N36").Value = -31880.334
$ws.Range("H107").Value = 10798.4
$ws.Range("I107").Value = 10798.4
$ws.Range("K107").Value = 32395.2
$ws.Range("M107").Value = -30475.2
$ws.Range("H122").Value = 2287.1853
$ws.Range("I122").Value = 1776.4445
$ws.Range("J122").Value = 3308.6667
$ws.Range("K122").Value = 5329.333500000001
$ws.Range("L122").Value = 9926.000100000001
$ws.Range("M122").Value = -2879.333500000001
$ws.Range("N122").Value = -14826.0001
